$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters (in order) corresponding to the values stored per data row
$cols = @("B","D","E","F","G","H","I","M","N")

# New values for rows 2..25 (each inner array holds values for $cols, in order)
$data = @(
    ,@(0.1530606109953538,0.01154117350315431,0.07351589593030106,1.074084204914811,1.063847557083534,0.7697743043384833,0.3362160722058078,0.9204868426998587,1.166626603982138)
    ,@(0.1357775162107089,0.01012326647565942,0.07544462030425381,0.9808143043797344,0.9535298163640391,0.7246785516775844,0.3385680073173276,0.8060981082078627,1.140237364641536)
    ,@(0.1251454038286255,0.009261683489590666,0.07669947880999484,0.9243249943486092,0.8864856452858021,0.6976036091732567,0.3402367490133109,0.7359037417032823,1.12446952433541)
    ,@(0.1208080467930444,0.008912713245010195,0.07722857663287064,0.9014962860631783,0.8593328219753857,0.6867222701876869,0.3409733869191101,0.7073075512665099,1.118154575476808)
    ,@(0.1200875608181349,0.008854891821606259,0.0773175032557889,0.8977170052329484,0.8548341132772634,0.6849245357964548,0.3410991279814652,0.7025596462919452,1.11711269477135)
    ,@(0.1250869271517985,0.009256968701372159,0.07670654263909027,0.924016350619425,0.8861187790852512,0.697456247766155,0.3402464541821182,0.7355180508983636,1.124383909510655)
    ,@(0.1471058310767859,0.01105032498063707,0.07416625000953347,1.041760238572977,1.025663066102538,0.7540963267523182,0.3369804675519816,0.8810360903169254,1.157437915459255)
    ,@(0.1901089266194163,0.01464497206745108,0.06974604841467702,1.279060346980458,1.305064320238671,0.8701589278946926,0.3323532324168532,1.166811393090484,1.225665764728546)
    ,@(0.221578419272646,0.01734250503623258,0.06684204630613699,1.45764155198205,1.514245240289199,0.9586495069499392,0.3300307661997302,1.377173024363415,1.277818279842307)
    ,@(0.2358640717861249,0.01858396967108433,0.06559583072891062,1.539876863569305,1.610343168017721,0.9996446275310973,0.3292068331565119,1.472993933997117,1.301973137387648)
    ,@(0.2412690071209909,0.01905630532586144,0.06513471100908408,1.571166514698405,1.64687510112617,1.015277959156037,0.3289281668162722,1.509299436028826,1.311180862379985)
    ,@(0.2401051746788596,0.01895447812247397,0.06523354071491294,1.564421017391822,1.639000887379893,1.011906130231921,0.3289867014581134,1.501479464890693,1.30919512691014)
    ,@(0.236308836288913,0.01862278360650294,0.06555767758086883,1.542448065839181,1.613345793964641,1.000928582318238,0.329183239528561,1.475980392690559,1.302729450026419)
    ,@(0.233982840829924,0.01841990496169643,0.0657576276520615,1.529008555915169,1.597649961137051,0.9942188477979812,0.3293079634769676,1.460364167337843,1.298776922772561)
    ,@(0.2206441790498559,0.01726167441371729,0.06692499812982522,1.452287822223724,1.507984491938828,0.9559855240865431,0.3300892831081654,1.370913617320767,1.276248270734811)
    ,@(0.2124533640376853,0.01655492222548105,0.06766033213070166,1.405481540543008,1.453223092927573,0.9327222845900565,0.3306280891867175,1.316072388364475,1.262537132569264)
    ,@(0.2077394272100719,0.01614975753520298,0.06809031798426002,1.378653394049735,1.421813911346504,0.9194113885291415,0.3309598927157857,1.284540930510573,1.25469142545461)
    ,@(0.2061429017249168,0.01601280121402482,0.06823711209332739,1.369585754576377,1.411194242081081,0.9149163973855821,0.3310759996354946,1.273866880759016,1.252042006396948)
    ,@(0.2133255827780829,0.01663001733761149,0.0675813255869846,1.410454423994253,1.459043369891788,0.9351914818313389,0.3305684669124851,1.321909098500853,1.263992515017264)
    ,@(0.2374240446401927,0.01872014880842698,0.06546217746994376,1.548897968758297,1.620877417836653,1.004149962484234,0.3291246075330534,1.4834695234306,1.304626934431013)
    ,@(0.2531460121449811,0.02009918523611987,0.06414013037422883,1.640249465548948,1.72747455454487,1.049857197097424,0.3283752417875654,1.58917735517565,1.331537799030713)
    ,@(0.2447575814097149,0.01936192392242475,0.06483996119360658,1.591411965510304,1.670503604834153,1.025402950268528,0.3287574495219516,1.532747539160823,1.317142935662361)
    ,@(0.2129312680036719,0.01659606324762564,0.06761702192029961,1.408205930227012,1.456411793058095,0.9340749603447875,0.3305953534953012,1.319270329802563,1.263334421255394)
    ,@(0.178496116742707,0.01366325940277591,0.07088160355633732,1.2141446965623,1.228822163494527,0.8382084943262953,0.3334154099205904,1.089441956060696,1.206848876835522)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $row).Value = $vals[$j]
    }
}
